$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the price list date (row 1) by one day
$ws.Range("A1").Value = 45311

# Update unit prices (column D) to the newly discounted values
$ws.Range("D14").Value = 4375.491
$ws.Range("D15").Value = 4874.98
$ws.Range("D16").Value = 6081.759
$ws.Range("D17").Value = 7392.386
$ws.Range("D18").Value = 8501.228999999999
$ws.Range("D19").Value = 9000.715
$ws.Range("D20").Value = 12997.743
$ws.Range("D21").Value = 11188.472
$ws.Range("D22").Value = 13126.483
$ws.Range("D23").Value = 14161.392
$ws.Range("D24").Value = 10462.154
$ws.Range("D25").Value = 9739.986000000001
$ws.Range("D26").Value = 11428.204
$ws.Range("D28").Value = 11687.961
$ws.Range("D29").Value = 14105.485
$ws.Range("D30").Value = 17501.958
$ws.Range("D31").Value = 20418.983
